# "adding averages and more checks"
$wb = $excel.ActiveWorkbook

$wsTraining = $wb.Worksheets.Item("Training Dashboard")
$wsExam = $wb.Worksheets.Item("Exam Dashboard")

# --- Header / title styling ----------------------------------------------
# Title (row 1) and the blue header band (row 2) both end up bold with a
# white font; the title no longer needs its larger 14pt size, so both
# rows converge on the same bold + white font.
$wsTraining.Range("A1").Font.Size = 11
$wsTraining.Range("A1").Font.Color = 16777215
$wsTraining.Range("A2:K2").Font.Color = 16777215

$wsExam.Range("A1").Font.Size = 11
$wsExam.Range("A1").Font.Color = 16777215
$wsExam.Range("A2:G2").Font.Color = 16777215

# --- Exam Dashboard sheet specific edits ---------------------------------
# Widen the COMMENTS column (target stored width of 15).
$wsExam.Columns.Item(5).ColumnWidth = 14.166666666666666

# Update the comment text for the first three data rows.
$wsExam.Range("E3").Value = "date is valid"
$wsExam.Range("E4").Value = "date is valid"
$wsExam.Range("E5").Value = "date is valid"
